$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "banding date" header in P1, matching the format of the other
# header cells (copy format from O1, the previous last header cell).
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P1").Value = "banding date"

# Style-1 (m/d/yyyy, numFmtId 14) template cell already used elsewhere in
# the sheet for the same kind of date value.
$ws.Range("C11").Copy()

$style1Rows = @(2, 6, 11, 13, 18, 22, 28, 33, 37, 38, 41, 46, 51, 54, 58, 63, 67, 72, 77, 82, 87, 93, 98)
foreach ($r in $style1Rows) {
    $ws.Range("P$r").PasteSpecial(-4122)   # xlPasteFormats
}

# Style-4 (d-mmm, numFmtId 16) template cell, used for P24.
$ws.Range("C2").Copy()
$ws.Range("P24").PasteSpecial(-4122)   # xlPasteFormats

# Now fill in the banding date values (Excel date serials).
$ws.Range("P2").Value = 44353
$ws.Range("P6").Value = 44355
$ws.Range("P11").Value = 44368
$ws.Range("P13").Value = 44362
$ws.Range("P18").Value = 44340
$ws.Range("P22").Value = 44362
$ws.Range("P24").Value = 44368
$ws.Range("P28").Value = 44362
$ws.Range("P33").Value = 44349
$ws.Range("P37").Value = 44368
$ws.Range("P38").Value = 44368
$ws.Range("P41").Value = 44340
$ws.Range("P46").Value = 44363
$ws.Range("P51").Value = 44353
$ws.Range("P54").Value = 44362
$ws.Range("P58").Value = 44349
$ws.Range("P63").Value = 44393
$ws.Range("P67").Value = 44355
$ws.Range("P72").Value = 44349
$ws.Range("P77").Value = 44355
$ws.Range("P82").Value = 44353
$ws.Range("P87").Value = 44349
$ws.Range("P93").Value = 44355
$ws.Range("P98").Value = 44393

# Match the final selection/viewport recorded in the saved workbook.
$ws.Range("P24").Select()
